$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "QfVkf357"
$ws.Range("B2").Value = 23080917
$ws.Range("C2").Value = "xeoloyz63"
$ws.Range("D2").Value = "Ps4u`$6%U"
$ws.Range("F2").Value = "kiWrVdoY"
$ws.Range("G2").Value = "FCgB"
